$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 769, shifting existing rows 769-810 down to 770-811
$ws.Rows.Item(769).Insert()

# Populate the newly inserted row with the new data point
# Force column A to be stored as plain text (matching the other date-like
# strings in this column) rather than letting Excel auto-convert it to a
# date serial number.
$ws.Cells.Item(769, 1).NumberFormat = "@"
$ws.Cells.Item(769, 1).Value = "2026/02/04"
$ws.Cells.Item(769, 2).Value = "水"
$ws.Cells.Item(769, 3).Value = 14
$ws.Cells.Item(769, 4).Value = 201
